# Aula 08 - Python Big Data
# Slide 2 title placeholder: bump the lesson number and swap the subtitle
# from "Contextualização" to "Data Science", keeping the existing run
# formatting (bold, white scheme color, 36pt for the subtitle line) intact.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

# --- "Aula 01" -> "Aula 08" -------------------------------------------------
$oldLesson = "Aula 01"
$newLesson = "Aula 08"
$full = $tr.Text
$pos = $full.IndexOf($oldLesson)
if ($pos -ge 0) {
    $lessonRange = $tr.Characters($pos + 1, $oldLesson.Length)
    $lessonRange.Text = $newLesson
}

# --- "Contextualização" -> "Data Science" -----------------------------------
$oldTopic = "Contextualização"
$newTopic = "Data Science"
$full = $tr.Text
$pos = $full.IndexOf($oldTopic)
if ($pos -ge 0) {
    $topicRange = $tr.Characters($pos + 1, $oldTopic.Length)
    $topicRange.Text = $newTopic
}

Write-Output "Title text now: $($tr.Text)"
